$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- New header row (row 27), columns H:M mirroring A:F headers ---
$ws1.Range("H27").Value = "Tevap"
$ws1.Range("I27").Value = "x1"
$ws1.Range("J27").Value = "x2"
$ws1.Range("K27").Value = "x3"
$ws1.Range("L27").Value = "pevap / Pa"
$ws1.Range("M27").Value = "in bar "

# --- Row 28: first data row of the new mini table ---
$ws1.Range("H28").Value = 25
$ws1.Range("I28").Formula = "=0.5"
$ws1.Range("J28").Value = 0.4
$ws1.Range("K28").Formula = "=1-I28-J28"
$ws1.Range("M28").Value = 6.5581

# --- Row 29 ---
$ws1.Range("I29").Formula = "=I28+0.05"
$ws1.Range("J29").Formula = "=J28-0.04"
$ws1.Range("K29").Formula = "=1-I29-J29"
$ws1.Range("M29").Value = 6.2548

# --- Rows 30-36: fill down the same formulas ---
$ws1.Range("I30:I36").Formula = "=I29+0.05"
$ws1.Range("J30:J36").Formula = "=J29-0.04"
$ws1.Range("K30:K36").Formula = "=1-I30-J30"

$ws1.Range("M30").Value = 5.9541
$ws1.Range("M31").Value = 5.655
$ws1.Range("M32").Value = 5.3563
$ws1.Range("M33").Value = 5.0569
$ws1.Range("M34").Value = 4.7556
$ws1.Range("M35").Value = 4.4512
$ws1.Range("M36").Value = 4.1425

# --- Column L (28-36): empty cells, formatted with a new style ---
# (10pt Arial, #,##0 number format, vertically centered)
# Seed L28 from an existing centered style (E28) so the new font/format
# combination is derived cleanly, then propagate the resulting format.
$seed = $ws1.Range("L28")
$ws1.Range("E28").Copy()
$seed.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$seed.Font.Name = "Arial"
$seed.Font.Size = 10
$seed.NumberFormat = "#,##0"

$seed.Copy()
$ws1.Range("L29:L36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- View / window settings ---
$ws2.Activate()
$excel.ActiveWindow.Zoom = 70

$ws1.Activate()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("V17").Select()
